# chore: adapt column header formatting to respective input file names
#
# - Renames the "_old" / "_new" header-name suffixes used throughout the
#   sheet's column headers (A1:U1) to "_FV2404" / "_FV2410" respectively
#   (the "diff" header in K1 is left untouched).
# - Wraps the data range A1:U79 in an Excel Table ("Table1") whose column
#   names mirror the new headers.
# - Freezes the header row (row 1) so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Rename the column headers in row 1.
# ---------------------------------------------------------------------
$headers = @(
  "Segmentname_FV2404", "Segmentgruppe_FV2404", "Segment_FV2404", "Datenelement_FV2404", "Segment ID_FV2404",
  "Code_FV2404", "Qualifier_FV2404", "Beschreibung_FV2404", "Bedingungsausdruck_FV2404", "Bedingung_FV2404",
  "diff",
  "Segmentname_FV2410", "Segmentgruppe_FV2410", "Segment_FV2410", "Datenelement_FV2410", "Segment ID_FV2410",
  "Code_FV2410", "Qualifier_FV2410", "Beschreibung_FV2410", "Bedingungsausdruck_FV2410", "Bedingung_FV2410"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# ---------------------------------------------------------------------
# 2. Turn A1:U79 into a proper Excel Table, keeping the existing header
#    row formatting (bold / grey fill / border / centered+wrapped text)
#    intact. Adding a ListObject in Excel always bakes the *current*
#    header formatting into a one-off "header row" dxf, so we stash the
#    formatting first, reset the header to the default style, add the
#    table (header dxf now captures the plain default look), and then
#    restore the original formatting on top - leaving the workbook's
#    style table untouched.
# ---------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$scratch = $ws.Range("A200:U200")

$headerRange.Copy() | Out-Null
$scratch.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$headerRange.Style = "Normal"

$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U79"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"
$lo.TableStyle = ""

$scratch.Copy() | Out-Null
$headerRange.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$scratch.Clear() | Out-Null

# ---------------------------------------------------------------------
# 3. Freeze the header row.
# ---------------------------------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

Write-Host "Header renaming, table creation and freeze pane applied."
